# Cryptos-list price/volume refresh (GitHub Actions scrape update).
# Rewrites the B:E inline-string cells that changed between runs -
# prices in column D, 1h volume deltas in column E, and for the two
# coins that swapped rank (rows 31/32) the name + link as well.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.865.53"
$ws.Range("E2").Value = "  +1.21%  "

# Row 3
$ws.Range("D3").Value = "3.119.29"
$ws.Range("E3").Value = "  +0.66%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").Value = "'599.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "

# Row 6
$ws.Range("D6").Value = "'141.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.08%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("D8").Value = "3.113.63"
$ws.Range("E8").Value = "  +0.37%  "

# Row 9
$ws.Range("E9").Value = "  +0.59%  "

# Row 10
$ws.Range("E10").Value = "  -0.33%  "

# Row 11
$ws.Range("D11").Value = "'5.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.29%  "

# Row 12
$ws.Range("E12").Value = "  +0.17%  "

# Row 13
$ws.Range("E13").Value = "  +1.07%  "

# Row 14
$ws.Range("D14").Value = "'34.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.31%  "

# Row 15
$ws.Range("D15").Value = "3.633.90"
$ws.Range("E15").Value = "  +0.73%  "

# Row 16
$ws.Range("E16").Value = "  +2.95%  "

# Row 17
$ws.Range("D17").Value = "63.939.45"
$ws.Range("E17").Value = "  +1.08%  "

# Row 18
$ws.Range("D18").Value = "3.125.83"
$ws.Range("E18").Value = "  +0.87%  "

# Row 19
$ws.Range("D19").Value = "'6.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.76%  "

# Row 20
$ws.Range("D20").Value = "'478.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.55%  "

# Row 21
$ws.Range("D21").Value = "'14.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "

# Row 22
$ws.Range("E22").Value = "  +0.81%  "

# Row 23
$ws.Range("D23").Value = "'7.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.68%  "

# Row 24
$ws.Range("D24").Value = "'87.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.97%  "

# Row 25
$ws.Range("D25").Value = "'13.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.05%  "

# Row 27
$ws.Range("D27").Value = "'2.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.53%  "

# Row 28
$ws.Range("D28").Value = "'8.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.32%  "

# Row 29
$ws.Range("D29").Value = "'7.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.38%  "

# Row 30
$ws.Range("E30").Value = "  -1.16%  "

# Row 31
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.16%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.111"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.08%  "

# Row 33
$ws.Range("D33").Value = "'26.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.10%  "

# Row 34
$ws.Range("E34").Value = "  -0.30%  "

# Row 35
$ws.Range("D35").Value = "'1.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.74%  "

# Row 36
$ws.Range("D36").Value = "'5.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.10%  "

# Row 37
$ws.Range("D37").Value = "'52.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0745"
$ws.Range("E38").Value = "  +0.47%  "

# Row 39
$ws.Range("E39").Value = "  +1.58%  "

# Row 40
$ws.Range("D40").Value = "'433.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.22%  "

# Row 41
$ws.Range("E41").Value = "  -0.07%  "

# Row 42
$ws.Range("E42").Value = "  +1.23%  "

# Row 43
$ws.Range("D43").Value = "'8.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.21%  "

# Row 44
$ws.Range("D44").Value = "2.860.76"
$ws.Range("E44").Value = "  +1.17%  "

# Row 45
$ws.Range("E45").Value = "  -2.27%  "

# Row 46
$ws.Range("E46").Value = "  -2.42%  "

# Row 47
$ws.Range("E47").Value = "  +0.39%  "

# Row 48
$ws.Range("E48").Value = "  +0.03%  "

# Row 49
$ws.Range("D49").Value = "'25.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.72%  "

# Row 50
$ws.Range("E50").Value = "  +0.33%  "

# Row 51
$ws.Range("D51").Value = "'121.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.37%  "

